# Add a "Save" column (H) to the s_vals sheet, mirroring the header style
# used by the existing columns (B:G) and filling row 2 with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it reuses the same cell style (bold, centered, bordered)
# instead of creating a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
